# Shift the daily/hourly forecast window forward (18.02-25.02 -> 20.02-27.02)
# and refresh the Prediction values with the retrained model output
# ("Updating the models with january production data").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$n = 169
$arr = New-Object 'object[,]' $n,4

$arr[0,0] = 46073
$arr[0,1] = 20
$arr[0,2] = 0.012
$arr[0,3] = "20.02.202620"
$arr[1,0] = 46073
$arr[1,1] = 21
$arr[1,2] = 0.012
$arr[1,3] = "20.02.202621"
$arr[2,0] = 46073
$arr[2,1] = 22
$arr[2,2] = 0.031
$arr[2,3] = "20.02.202622"
$arr[3,0] = 46073
$arr[3,1] = 23
$arr[3,2] = 0.031
$arr[3,3] = "20.02.202623"
$arr[4,0] = 46073
$arr[4,1] = 24
$arr[4,2] = 0.031
$arr[4,3] = "20.02.202624"
$arr[5,0] = 46074
$arr[5,1] = 1
$arr[5,2] = 0.032
$arr[5,3] = "21.02.20261"
$arr[6,0] = 46074
$arr[6,1] = 2
$arr[6,2] = 0.032
$arr[6,3] = "21.02.20262"
$arr[7,0] = 46074
$arr[7,1] = 3
$arr[7,2] = 0.032
$arr[7,3] = "21.02.20263"
$arr[8,0] = 46074
$arr[8,1] = 4
$arr[8,2] = 0.032
$arr[8,3] = "21.02.20264"
$arr[9,0] = 46074
$arr[9,1] = 5
$arr[9,2] = 0.032
$arr[9,3] = "21.02.20265"
$arr[10,0] = 46074
$arr[10,1] = 6
$arr[10,2] = 0.032
$arr[10,3] = "21.02.20266"
$arr[11,0] = 46074
$arr[11,1] = 7
$arr[11,2] = 0.032
$arr[11,3] = "21.02.20267"
$arr[12,0] = 46074
$arr[12,1] = 8
$arr[12,2] = 0.032
$arr[12,3] = "21.02.20268"
$arr[13,0] = 46074
$arr[13,1] = 9
$arr[13,2] = 0.051
$arr[13,3] = "21.02.20269"
$arr[14,0] = 46074
$arr[14,1] = 10
$arr[14,2] = 0.073
$arr[14,3] = "21.02.202610"
$arr[15,0] = 46074
$arr[15,1] = 11
$arr[15,2] = 0.292
$arr[15,3] = "21.02.202611"
$arr[16,0] = 46074
$arr[16,1] = 12
$arr[16,2] = 0.462
$arr[16,3] = "21.02.202612"
$arr[17,0] = 46074
$arr[17,1] = 13
$arr[17,2] = 0.596
$arr[17,3] = "21.02.202613"
$arr[18,0] = 46074
$arr[18,1] = 14
$arr[18,2] = 0.619
$arr[18,3] = "21.02.202614"
$arr[19,0] = 46074
$arr[19,1] = 15
$arr[19,2] = 0.572
$arr[19,3] = "21.02.202615"
$arr[20,0] = 46074
$arr[20,1] = 16
$arr[20,2] = 0.425
$arr[20,3] = "21.02.202616"
$arr[21,0] = 46074
$arr[21,1] = 17
$arr[21,2] = 0.282
$arr[21,3] = "21.02.202617"
$arr[22,0] = 46074
$arr[22,1] = 18
$arr[22,2] = 0.071
$arr[22,3] = "21.02.202618"
$arr[23,0] = 46074
$arr[23,1] = 19
$arr[23,2] = 0.045
$arr[23,3] = "21.02.202619"
$arr[24,0] = 46074
$arr[24,1] = 20
$arr[24,2] = 0.031
$arr[24,3] = "21.02.202620"
$arr[25,0] = 46074
$arr[25,1] = 21
$arr[25,2] = 0.031
$arr[25,3] = "21.02.202621"
$arr[26,0] = 46074
$arr[26,1] = 22
$arr[26,2] = 0.031
$arr[26,3] = "21.02.202622"
$arr[27,0] = 46074
$arr[27,1] = 23
$arr[27,2] = 0.031
$arr[27,3] = "21.02.202623"
$arr[28,0] = 46074
$arr[28,1] = 24
$arr[28,2] = 0.031
$arr[28,3] = "21.02.202624"
$arr[29,0] = 46075
$arr[29,1] = 1
$arr[29,2] = 0.032
$arr[29,3] = "22.02.20261"
$arr[30,0] = 46075
$arr[30,1] = 2
$arr[30,2] = 0.032
$arr[30,3] = "22.02.20262"
$arr[31,0] = 46075
$arr[31,1] = 3
$arr[31,2] = 0.032
$arr[31,3] = "22.02.20263"
$arr[32,0] = 46075
$arr[32,1] = 4
$arr[32,2] = 0.032
$arr[32,3] = "22.02.20264"
$arr[33,0] = 46075
$arr[33,1] = 5
$arr[33,2] = 0.032
$arr[33,3] = "22.02.20265"
$arr[34,0] = 46075
$arr[34,1] = 6
$arr[34,2] = 0.032
$arr[34,3] = "22.02.20266"
$arr[35,0] = 46075
$arr[35,1] = 7
$arr[35,2] = 0.032
$arr[35,3] = "22.02.20267"
$arr[36,0] = 46075
$arr[36,1] = 8
$arr[36,2] = 0.032
$arr[36,3] = "22.02.20268"
$arr[37,0] = 46075
$arr[37,1] = 9
$arr[37,2] = 0.051
$arr[37,3] = "22.02.20269"
$arr[38,0] = 46075
$arr[38,1] = 10
$arr[38,2] = 0.152
$arr[38,3] = "22.02.202610"
$arr[39,0] = 46075
$arr[39,1] = 11
$arr[39,2] = 0.348
$arr[39,3] = "22.02.202611"
$arr[40,0] = 46075
$arr[40,1] = 12
$arr[40,2] = 0.462
$arr[40,3] = "22.02.202612"
$arr[41,0] = 46075
$arr[41,1] = 13
$arr[41,2] = 0.585
$arr[41,3] = "22.02.202613"
$arr[42,0] = 46075
$arr[42,1] = 14
$arr[42,2] = 0.569
$arr[42,3] = "22.02.202614"
$arr[43,0] = 46075
$arr[43,1] = 15
$arr[43,2] = 0.544
$arr[43,3] = "22.02.202615"
$arr[44,0] = 46075
$arr[44,1] = 16
$arr[44,2] = 0.432
$arr[44,3] = "22.02.202616"
$arr[45,0] = 46075
$arr[45,1] = 17
$arr[45,2] = 0.288
$arr[45,3] = "22.02.202617"
$arr[46,0] = 46075
$arr[46,1] = 18
$arr[46,2] = 0.144
$arr[46,3] = "22.02.202618"
$arr[47,0] = 46075
$arr[47,1] = 19
$arr[47,2] = 0.045
$arr[47,3] = "22.02.202619"
$arr[48,0] = 46075
$arr[48,1] = 20
$arr[48,2] = 0.031
$arr[48,3] = "22.02.202620"
$arr[49,0] = 46075
$arr[49,1] = 21
$arr[49,2] = 0.031
$arr[49,3] = "22.02.202621"
$arr[50,0] = 46075
$arr[50,1] = 22
$arr[50,2] = 0.031
$arr[50,3] = "22.02.202622"
$arr[51,0] = 46075
$arr[51,1] = 23
$arr[51,2] = 0.031
$arr[51,3] = "22.02.202623"
$arr[52,0] = 46075
$arr[52,1] = 24
$arr[52,2] = 0.031
$arr[52,3] = "22.02.202624"
$arr[53,0] = 46076
$arr[53,1] = 1
$arr[53,2] = 0.032
$arr[53,3] = "23.02.20261"
$arr[54,0] = 46076
$arr[54,1] = 2
$arr[54,2] = 0.032
$arr[54,3] = "23.02.20262"
$arr[55,0] = 46076
$arr[55,1] = 3
$arr[55,2] = 0.032
$arr[55,3] = "23.02.20263"
$arr[56,0] = 46076
$arr[56,1] = 4
$arr[56,2] = 0.032
$arr[56,3] = "23.02.20264"
$arr[57,0] = 46076
$arr[57,1] = 5
$arr[57,2] = 0.032
$arr[57,3] = "23.02.20265"
$arr[58,0] = 46076
$arr[58,1] = 6
$arr[58,2] = 0.032
$arr[58,3] = "23.02.20266"
$arr[59,0] = 46076
$arr[59,1] = 7
$arr[59,2] = 0.032
$arr[59,3] = "23.02.20267"
$arr[60,0] = 46076
$arr[60,1] = 8
$arr[60,2] = 0.032
$arr[60,3] = "23.02.20268"
$arr[61,0] = 46076
$arr[61,1] = 9
$arr[61,2] = 0.051
$arr[61,3] = "23.02.20269"
$arr[62,0] = 46076
$arr[62,1] = 10
$arr[62,2] = 0.059
$arr[62,3] = "23.02.202610"
$arr[63,0] = 46076
$arr[63,1] = 11
$arr[63,2] = 0.279
$arr[63,3] = "23.02.202611"
$arr[64,0] = 46076
$arr[64,1] = 12
$arr[64,2] = 0.334
$arr[64,3] = "23.02.202612"
$arr[65,0] = 46076
$arr[65,1] = 13
$arr[65,2] = 0.334
$arr[65,3] = "23.02.202613"
$arr[66,0] = 46076
$arr[66,1] = 14
$arr[66,2] = 0.334
$arr[66,3] = "23.02.202614"
$arr[67,0] = 46076
$arr[67,1] = 15
$arr[67,2] = 0.28
$arr[67,3] = "23.02.202615"
$arr[68,0] = 46076
$arr[68,1] = 16
$arr[68,2] = 0.274
$arr[68,3] = "23.02.202616"
$arr[69,0] = 46076
$arr[69,1] = 17
$arr[69,2] = 0.182
$arr[69,3] = "23.02.202617"
$arr[70,0] = 46076
$arr[70,1] = 18
$arr[70,2] = 0.037
$arr[70,3] = "23.02.202618"
$arr[71,0] = 46076
$arr[71,1] = 19
$arr[71,2] = 0.026
$arr[71,3] = "23.02.202619"
$arr[72,0] = 46076
$arr[72,1] = 20
$arr[72,2] = 0.013
$arr[72,3] = "23.02.202620"
$arr[73,0] = 46076
$arr[73,1] = 21
$arr[73,2] = 0.013
$arr[73,3] = "23.02.202621"
$arr[74,0] = 46076
$arr[74,1] = 22
$arr[74,2] = 0.013
$arr[74,3] = "23.02.202622"
$arr[75,0] = 46076
$arr[75,1] = 23
$arr[75,2] = 0.013
$arr[75,3] = "23.02.202623"
$arr[76,0] = 46076
$arr[76,1] = 24
$arr[76,2] = 0.013
$arr[76,3] = "23.02.202624"
$arr[77,0] = 46077
$arr[77,1] = 1
$arr[77,2] = 0.014
$arr[77,3] = "24.02.20261"
$arr[78,0] = 46077
$arr[78,1] = 2
$arr[78,2] = 0.014
$arr[78,3] = "24.02.20262"
$arr[79,0] = 46077
$arr[79,1] = 3
$arr[79,2] = 0.032
$arr[79,3] = "24.02.20263"
$arr[80,0] = 46077
$arr[80,1] = 4
$arr[80,2] = 0.032
$arr[80,3] = "24.02.20264"
$arr[81,0] = 46077
$arr[81,1] = 5
$arr[81,2] = 0.032
$arr[81,3] = "24.02.20265"
$arr[82,0] = 46077
$arr[82,1] = 6
$arr[82,2] = 0.032
$arr[82,3] = "24.02.20266"
$arr[83,0] = 46077
$arr[83,1] = 7
$arr[83,2] = 0.032
$arr[83,3] = "24.02.20267"
$arr[84,0] = 46077
$arr[84,1] = 8
$arr[84,2] = 0.032
$arr[84,3] = "24.02.20268"
$arr[85,0] = 46077
$arr[85,1] = 9
$arr[85,2] = 0.038
$arr[85,3] = "24.02.20269"
$arr[86,0] = 46077
$arr[86,1] = 10
$arr[86,2] = 0.139
$arr[86,3] = "24.02.202610"
$arr[87,0] = 46077
$arr[87,1] = 11
$arr[87,2] = 0.334
$arr[87,3] = "24.02.202611"
$arr[88,0] = 46077
$arr[88,1] = 12
$arr[88,2] = 0.561
$arr[88,3] = "24.02.202612"
$arr[89,0] = 46077
$arr[89,1] = 13
$arr[89,2] = 0.742
$arr[89,3] = "24.02.202613"
$arr[90,0] = 46077
$arr[90,1] = 14
$arr[90,2] = 0.728
$arr[90,3] = "24.02.202614"
$arr[91,0] = 46077
$arr[91,1] = 15
$arr[91,2] = 0.489
$arr[91,3] = "24.02.202615"
$arr[92,0] = 46077
$arr[92,1] = 16
$arr[92,2] = 0.39
$arr[92,3] = "24.02.202616"
$arr[93,0] = 46077
$arr[93,1] = 17
$arr[93,2] = 0.285
$arr[93,3] = "24.02.202617"
$arr[94,0] = 46077
$arr[94,1] = 18
$arr[94,2] = 0.125
$arr[94,3] = "24.02.202618"
$arr[95,0] = 46077
$arr[95,1] = 19
$arr[95,2] = 0.026
$arr[95,3] = "24.02.202619"
$arr[96,0] = 46077
$arr[96,1] = 20
$arr[96,2] = 0.013
$arr[96,3] = "24.02.202620"
$arr[97,0] = 46077
$arr[97,1] = 21
$arr[97,2] = 0.031
$arr[97,3] = "24.02.202621"
$arr[98,0] = 46077
$arr[98,1] = 22
$arr[98,2] = 0.031
$arr[98,3] = "24.02.202622"
$arr[99,0] = 46077
$arr[99,1] = 23
$arr[99,2] = 0.031
$arr[99,3] = "24.02.202623"
$arr[100,0] = 46077
$arr[100,1] = 24
$arr[100,2] = 0.031
$arr[100,3] = "24.02.202624"
$arr[101,0] = 46078
$arr[101,1] = 1
$arr[101,2] = 0.032
$arr[101,3] = "25.02.20261"
$arr[102,0] = 46078
$arr[102,1] = 2
$arr[102,2] = 0.032
$arr[102,3] = "25.02.20262"
$arr[103,0] = 46078
$arr[103,1] = 3
$arr[103,2] = 0.032
$arr[103,3] = "25.02.20263"
$arr[104,0] = 46078
$arr[104,1] = 4
$arr[104,2] = 0.032
$arr[104,3] = "25.02.20264"
$arr[105,0] = 46078
$arr[105,1] = 5
$arr[105,2] = 0.032
$arr[105,3] = "25.02.20265"
$arr[106,0] = 46078
$arr[106,1] = 6
$arr[106,2] = 0.032
$arr[106,3] = "25.02.20266"
$arr[107,0] = 46078
$arr[107,1] = 7
$arr[107,2] = 0.032
$arr[107,3] = "25.02.20267"
$arr[108,0] = 46078
$arr[108,1] = 8
$arr[108,2] = 0.032
$arr[108,3] = "25.02.20268"
$arr[109,0] = 46078
$arr[109,1] = 9
$arr[109,2] = 0.051
$arr[109,3] = "25.02.20269"
$arr[110,0] = 46078
$arr[110,1] = 10
$arr[110,2] = 0.183
$arr[110,3] = "25.02.202610"
$arr[111,0] = 46078
$arr[111,1] = 11
$arr[111,2] = 0.348
$arr[111,3] = "25.02.202611"
$arr[112,0] = 46078
$arr[112,1] = 12
$arr[112,2] = 0.585
$arr[112,3] = "25.02.202612"
$arr[113,0] = 46078
$arr[113,1] = 13
$arr[113,2] = 0.805
$arr[113,3] = "25.02.202613"
$arr[114,0] = 46078
$arr[114,1] = 14
$arr[114,2] = 0.858
$arr[114,3] = "25.02.202614"
$arr[115,0] = 46078
$arr[115,1] = 15
$arr[115,2] = 0.863
$arr[115,3] = "25.02.202615"
$arr[116,0] = 46078
$arr[116,1] = 16
$arr[116,2] = 0.795
$arr[116,3] = "25.02.202616"
$arr[117,0] = 46078
$arr[117,1] = 17
$arr[117,2] = 0.582
$arr[117,3] = "25.02.202617"
$arr[118,0] = 46078
$arr[118,1] = 18
$arr[118,2] = 0.297
$arr[118,3] = "25.02.202618"
$arr[119,0] = 46078
$arr[119,1] = 19
$arr[119,2] = 0.045
$arr[119,3] = "25.02.202619"
$arr[120,0] = 46078
$arr[120,1] = 20
$arr[120,2] = 0.031
$arr[120,3] = "25.02.202620"
$arr[121,0] = 46078
$arr[121,1] = 21
$arr[121,2] = 0.031
$arr[121,3] = "25.02.202621"
$arr[122,0] = 46078
$arr[122,1] = 22
$arr[122,2] = 0.031
$arr[122,3] = "25.02.202622"
$arr[123,0] = 46078
$arr[123,1] = 23
$arr[123,2] = 0.031
$arr[123,3] = "25.02.202623"
$arr[124,0] = 46078
$arr[124,1] = 24
$arr[124,2] = 0.031
$arr[124,3] = "25.02.202624"
$arr[125,0] = 46079
$arr[125,1] = 1
$arr[125,2] = 0.032
$arr[125,3] = "26.02.20261"
$arr[126,0] = 46079
$arr[126,1] = 2
$arr[126,2] = 0.032
$arr[126,3] = "26.02.20262"
$arr[127,0] = 46079
$arr[127,1] = 3
$arr[127,2] = 0.032
$arr[127,3] = "26.02.20263"
$arr[128,0] = 46079
$arr[128,1] = 4
$arr[128,2] = 0.032
$arr[128,3] = "26.02.20264"
$arr[129,0] = 46079
$arr[129,1] = 5
$arr[129,2] = 0.032
$arr[129,3] = "26.02.20265"
$arr[130,0] = 46079
$arr[130,1] = 6
$arr[130,2] = 0.032
$arr[130,3] = "26.02.20266"
$arr[131,0] = 46079
$arr[131,1] = 7
$arr[131,2] = 0.032
$arr[131,3] = "26.02.20267"
$arr[132,0] = 46079
$arr[132,1] = 8
$arr[132,2] = 0.032
$arr[132,3] = "26.02.20268"
$arr[133,0] = 46079
$arr[133,1] = 9
$arr[133,2] = 0.051
$arr[133,3] = "26.02.20269"
$arr[134,0] = 46079
$arr[134,1] = 10
$arr[134,2] = 0.294
$arr[134,3] = "26.02.202610"
$arr[135,0] = 46079
$arr[135,1] = 11
$arr[135,2] = 0.608
$arr[135,3] = "26.02.202611"
$arr[136,0] = 46079
$arr[136,1] = 12
$arr[136,2] = 0.886
$arr[136,3] = "26.02.202612"
$arr[137,0] = 46079
$arr[137,1] = 13
$arr[137,2] = 1.049
$arr[137,3] = "26.02.202613"
$arr[138,0] = 46079
$arr[138,1] = 14
$arr[138,2] = 1.111
$arr[138,3] = "26.02.202614"
$arr[139,0] = 46079
$arr[139,1] = 15
$arr[139,2] = 1.104
$arr[139,3] = "26.02.202615"
$arr[140,0] = 46079
$arr[140,1] = 16
$arr[140,2] = 0.877
$arr[140,3] = "26.02.202616"
$arr[141,0] = 46079
$arr[141,1] = 17
$arr[141,2] = 0.802
$arr[141,3] = "26.02.202617"
$arr[142,0] = 46079
$arr[142,1] = 18
$arr[142,2] = 0.297
$arr[142,3] = "26.02.202618"
$arr[143,0] = 46079
$arr[143,1] = 19
$arr[143,2] = 0.053
$arr[143,3] = "26.02.202619"
$arr[144,0] = 46079
$arr[144,1] = 20
$arr[144,2] = 0.031
$arr[144,3] = "26.02.202620"
$arr[145,0] = 46079
$arr[145,1] = 21
$arr[145,2] = 0.031
$arr[145,3] = "26.02.202621"
$arr[146,0] = 46079
$arr[146,1] = 22
$arr[146,2] = 0.031
$arr[146,3] = "26.02.202622"
$arr[147,0] = 46079
$arr[147,1] = 23
$arr[147,2] = 0.031
$arr[147,3] = "26.02.202623"
$arr[148,0] = 46079
$arr[148,1] = 24
$arr[148,2] = 0.031
$arr[148,3] = "26.02.202624"
$arr[149,0] = 46080
$arr[149,1] = 1
$arr[149,2] = 0.032
$arr[149,3] = "27.02.20261"
$arr[150,0] = 46080
$arr[150,1] = 2
$arr[150,2] = 0.032
$arr[150,3] = "27.02.20262"
$arr[151,0] = 46080
$arr[151,1] = 3
$arr[151,2] = 0.032
$arr[151,3] = "27.02.20263"
$arr[152,0] = 46080
$arr[152,1] = 4
$arr[152,2] = 0.032
$arr[152,3] = "27.02.20264"
$arr[153,0] = 46080
$arr[153,1] = 5
$arr[153,2] = 0.032
$arr[153,3] = "27.02.20265"
$arr[154,0] = 46080
$arr[154,1] = 6
$arr[154,2] = 0.032
$arr[154,3] = "27.02.20266"
$arr[155,0] = 46080
$arr[155,1] = 7
$arr[155,2] = 0.032
$arr[155,3] = "27.02.20267"
$arr[156,0] = 46080
$arr[156,1] = 8
$arr[156,2] = 0.032
$arr[156,3] = "27.02.20268"
$arr[157,0] = 46080
$arr[157,1] = 9
$arr[157,2] = 0.051
$arr[157,3] = "27.02.20269"
$arr[158,0] = 46080
$arr[158,1] = 10
$arr[158,2] = 0.35
$arr[158,3] = "27.02.202610"
$arr[159,0] = 46080
$arr[159,1] = 11
$arr[159,2] = 0.87
$arr[159,3] = "27.02.202611"
$arr[160,0] = 46080
$arr[160,1] = 12
$arr[160,2] = 0.953
$arr[160,3] = "27.02.202612"
$arr[161,0] = 46080
$arr[161,1] = 13
$arr[161,2] = 1.301
$arr[161,3] = "27.02.202613"
$arr[162,0] = 46080
$arr[162,1] = 14
$arr[162,2] = 1.291
$arr[162,3] = "27.02.202614"
$arr[163,0] = 46080
$arr[163,1] = 15
$arr[163,2] = 1.257
$arr[163,3] = "27.02.202615"
$arr[164,0] = 46080
$arr[164,1] = 16
$arr[164,2] = 1.256
$arr[164,3] = "27.02.202616"
$arr[165,0] = 46080
$arr[165,1] = 17
$arr[165,2] = 0.824
$arr[165,3] = "27.02.202617"
$arr[166,0] = 46080
$arr[166,1] = 18
$arr[166,2] = 0.34
$arr[166,3] = "27.02.202618"
$arr[167,0] = 46080
$arr[167,1] = 19
$arr[167,2] = 0.053
$arr[167,3] = "27.02.202619"
$arr[168,0] = 46080
$arr[168,1] = 20
$arr[168,2] = 0.031
$arr[168,3] = "27.02.202620"

# Write Data / Interval / Prediction / Lookup columns for rows 2-170 in one shot
$ws.Range("A2:D170").Value = $arr
